$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2359592616558075
$ws.Range("B1").Value = 0.2177231460809708
$ws.Range("C1").Value = 0.2155307680368423
$ws.Range("D1").Value = 0.2698128223419189
$ws.Range("E1").Value = 0.4104617834091187
